# Apply updates to the "ODI" column (column C) for several players whose
# match-count values increased by 1 as additional data was scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Active Players")

$rows = @(2, 8, 9, 14, 21, 23, 29, 30, 33, 39, 42)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 1
}
